$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-existing (but blank) row 17 needs to end up fully styled
# like the rest of column A once it gets a value. Cycling it out and back in
# makes sure it picks up the column's default style when content is written.
$ws.Rows("17:17").Delete()
$ws.Rows("17:17").Insert()

# Final list of file names for column A: the original nine entries plus
# "cryptoAixbtBase.xlsx" inserted alphabetically, followed by the eight
# newly uploaded files appended at the end.
$items = @(
    "cryptoAAVEpolygon.xlsx",
    "cryptoAerodromeBase.xlsx",
    "cryptoAixbtBase.xlsx",
    "cryptoChainlinkPolygon.xlsx",
    "cryptoDogeBnb.xlsx",
    "cryptoMorphoBase.xlsx",
    "cryptoMystPolygon.xlsx",
    "cryptoPaxgoldPolygon.xlsx",
    "cryptoSolanaPolygon.xlsx",
    "cryptoWrappedBTCPolygon.xlsx",
    "cryptoUniswapPolygon.xlsx",
    "cryptoAsterBnb.xlsx",
    "cryptoAtomcosmosBnb.xlsx",
    "cryptoNexoPolygon.xlsx",
    "cryptoVirtualBase.xlsx",
    "cryptoPancakeswapBase.xlsx",
    "cryptoWrappedEthPolygon.xlsx",
    "cryptoWrappedbnbBnb.xlsx"
)

$row = 1
foreach ($item in $items) {
    $ws.Range("A$row").Value = $item
    $row = $row + 1
}

$ws.Range("F23").Select()
